# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Reorders the "Periodo Mora" rows (16-31) into chronological order
# (2104 .. 2206) and inserts the JOHN JAIRO RUIZ ARIZA / 1047408026 entry
# right after the first INIRIDA ... 2104 row (now row 17) instead of at
# the very end of the table (row 31).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B = Tipo Doc Trabajador, C = N Doc Trabajador, D = Nombre Trabajador,
# E = Periodo Mora, F = Valor Mora, G = Salario Basico
$rows = @(
    @{ Row = 16; B = "CC"; C = "45752903";   D = "INIRIDA DEL CARMEN CAMARGO PEREZ"; E = "2104"; F = 27861; G = 908526 },
    @{ Row = 17; B = "CC"; C = "1047408026"; D = "JOHN JAIRO RUIZ ARIZA";            E = "2104"; F = 24000; G = 1500000 },
    @{ Row = 18; B = "CC"; C = "45752903";   D = "INIRIDA DEL CARMEN CAMARGO PEREZ"; E = "2105"; F = 36341; G = 908526 },
    @{ Row = 19; B = "CC"; C = "45752903";   D = "INIRIDA DEL CARMEN CAMARGO PEREZ"; E = "2106"; F = 36341; G = 908526 },
    @{ Row = 20; B = "CC"; C = "45752903";   D = "INIRIDA DEL CARMEN CAMARGO PEREZ"; E = "2107"; F = 36341; G = 908526 },
    @{ Row = 21; B = "CC"; C = "45752903";   D = "INIRIDA DEL CARMEN CAMARGO PEREZ"; E = "2108"; F = 36341; G = 908526 },
    @{ Row = 22; B = "CC"; C = "45752903";   D = "INIRIDA DEL CARMEN CAMARGO PEREZ"; E = "2109"; F = 36341; G = 908526 },
    @{ Row = 23; B = "CC"; C = "45752903";   D = "INIRIDA DEL CARMEN CAMARGO PEREZ"; E = "2110"; F = 36341; G = 908526 },
    @{ Row = 24; B = "CC"; C = "45752903";   D = "INIRIDA DEL CARMEN CAMARGO PEREZ"; E = "2111"; F = 36341; G = 908526 },
    @{ Row = 25; B = "CC"; C = "45752903";   D = "INIRIDA DEL CARMEN CAMARGO PEREZ"; E = "2112"; F = 36341; G = 908526 },
    @{ Row = 26; B = "CC"; C = "45752903";   D = "INIRIDA DEL CARMEN CAMARGO PEREZ"; E = "2201"; F = 36341; G = 908526 },
    @{ Row = 27; B = "CC"; C = "45752903";   D = "INIRIDA DEL CARMEN CAMARGO PEREZ"; E = "2202"; F = 36341; G = 908526 },
    @{ Row = 28; B = "CC"; C = "45752903";   D = "INIRIDA DEL CARMEN CAMARGO PEREZ"; E = "2203"; F = 36341; G = 908526 },
    @{ Row = 29; B = "CC"; C = "45752903";   D = "INIRIDA DEL CARMEN CAMARGO PEREZ"; E = "2204"; F = 36341; G = 908526 },
    @{ Row = 30; B = "CC"; C = "45752903";   D = "INIRIDA DEL CARMEN CAMARGO PEREZ"; E = "2205"; F = 36341; G = 908526 },
    @{ Row = 31; B = "CC"; C = "45752903";   D = "INIRIDA DEL CARMEN CAMARGO PEREZ"; E = "2206"; F = 26650; G = 908526 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("B$n").Value = $r.B
    $ws.Range("C$n").Value = $r.C
    $ws.Range("D$n").Value = $r.D
    $ws.Range("E$n").Value = $r.E
    $ws.Range("F$n").Value = $r.F
    $ws.Range("G$n").Value = $r.G
}
